$d = $word.ActiveDocument

$replacements = @(
    @{old = "73×69="; new = "55×64="},
    @{old = "66×27="; new = "51×44="},
    @{old = "13×28="; new = "62×43="},
    @{old = "12×49="; new = "61×21="},
    @{old = "45×84="; new = "99×99="},
    @{old = "58×27="; new = "33×72="},
    @{old = "20×16="; new = "32×95="},
    @{old = "29×97="; new = "95×33="},
    @{old = "76×20="; new = "92×42="},
    @{old = "34×19="; new = "21×66="},
    @{old = "19×14="; new = "63×22="},
    @{old = "68×18="; new = "99×13="},
    @{old = "95×76="; new = "46×11="},
    @{old = "38×30="; new = "80×28="},
    @{old = "25×44="; new = "13×44="},
    @{old = "14×27="; new = "20×45="},
    @{old = "70×30="; new = "16×28="},
    @{old = "25×67="; new = "24×41="},
    @{old = "70×31="; new = "75×12="},
    @{old = "90×57="; new = "46×58="},
    @{old = "20×98="; new = "24×25="},
    @{old = "30×93="; new = "49×45="},
    @{old = "99×89="; new = "62×38="},
    @{old = "38×68="; new = "87×86="},
    @{old = "97×14="; new = "43×25="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
